# Apply the "hydropower disaggregation, base case fix" edit.
#
# Adds two new attribute rows (NobjHYD / NobjIRR) describing the number of
# disaggregated objective variables the simulation outputs for hydropower
# production and irrigation, extends Table1 to cover the new rows, and
# updates the active selection on the ModelParameters sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ModelParameters")

# --- New row 20: NobjHYD ---------------------------------------------------
$ws.Range("B20").Value = "NobjHYD"
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = "int"
$ws.Range("E20").Value = "Number of objective variables that the simulation outputs for the disaggregated objectives for hydropower production"

# --- New row 21: NobjIRR ----------------------------------------------------
$ws.Range("B21").Value = "NobjIRR"
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = "int"
$ws.Range("E21").Value = "Number of objective variables that the simulation outputs for the disaggregated objectives for irrigation"

# --- Formatting for the two new rows (wrap text, vertically centred, -------
# --- Calibri 12) so they get their own cell style, matching the style -------
# --- used for the rest of the settings tables. ------------------------------
$newRows = $ws.Range("B20:E21")
$newRows.Font.Name = "Calibri"
$newRows.WrapText = $true
$newRows.VerticalAlignment = -4108

# --- Extend Table1 (the ListObject) so it covers the new rows. -------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B4:E21"))

# --- Update the active selection / cursor position on the sheet. ----------
$ws.Range("B24").Select() | Out-Null
